$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Request Header (E4) None -> ContentType.JSON
$ws.Range("E4").Value = "ContentType.JSON"

# Row 6: Query/header Parameters (I6) {"page": 2} -> page=2
#        Assertions (N6) None -> data: contains only one page with 6 users
$ws.Range("I6").Value = "page=2"
$ws.Range("N6").Value = "data: contains only one page with 6 users"

# Row 7: was /users/3 GET ... now becomes /register POST ...
$ws.Range("C7").Value = "/register"
$ws.Range("D7").Value = "POST"
$ws.Range("H7").Value = "None"
$ws.Range("J7").Value = "BasicAuth"
$ws.Range("K7").Value = "username: 'testuser', password: 'testpass'"
# Status code column is stored as text in this sheet, not a number - force text
$ws.Range("L7").Value = "'400"
$ws.Range("N7").Value = "None"

# Row 8: was /register POST ... now becomes /users/3 GET ...
$ws.Range("C8").Value = "/users/3"
$ws.Range("D8").Value = "GET"
$ws.Range("H8").Value = "id=3"
$ws.Range("J8").Value = "None"
$ws.Range("K8").Value = "None"
# Status code column is stored as text in this sheet, not a number - force text
$ws.Range("L8").Value = "'200"
$ws.Range("N8").Value = "data.last_name: 'Wong'"

# Row 9 removed entirely - delete the whole row, shifting nothing below it up
$ws.Range("A9:N9").EntireRow.Delete()
